$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Mesa"
$ws.Range("B1").Value = "Circuito"
$ws.Range("C1").Value = "Cantidad de Votantes"
$ws.Range("D1").Value = "FRENTE DE IZQUIERDA Y DE TRABAJADORES - UNIDAD"
$ws.Range("E1").Value = "CONSENSO FEDERAL"
$ws.Range("F1").Value = "JUNTOS POR EL CAMBIO"
$ws.Range("G1").Value = "FRENTE DE TODOS"
$ws.Range("H1").Value = "UNITE POR LA LIBERTAD Y LA DIGNIDAD"
$ws.Range("I1").Value = "Votos Nulos"
$ws.Range("J1").Value = "Votos Recurridos"
$ws.Range("K1").Value = "Votos impugnados"
$ws.Range("L1").Value = "Votos en blanco"

# --- Data row (row 2) ---
# A2/B2 hold identifiers that must stay textual (e.g. leading zeros in "00039"),
# so force text format while assigning, then restore the default style so the
# cell keeps using the workbook's normal style (same as every other cell).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1244"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "00039"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = 269
$ws.Range("D2").Value = 16
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 111
$ws.Range("G2").Value = 110
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 13
